$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("20-10-2025 07:30:49", "Nguyễn Văn Nam", "Đăng nhập."),
    @("20-10-2025 07:31:06", "Nguyễn Văn Nam", "Thêm nhân viên: Phan Thúy (Staff)"),
    @("20-10-2025 07:31:11", "Nguyễn Văn Nam", "Làm mới nhật ký hoạt động"),
    @("20-10-2025 07:31:14", "Nguyễn Văn Nam", "Đăng xuất"),
    @("20-10-2025 07:38:38", "Nguyễn Văn Nam", "Đăng nhập."),
    @("20-10-2025 07:38:53", "Nguyễn Văn Nam", "Thêm nhân viên: Phan Thúy (Staff)"),
    @("20-10-2025 07:38:55", "Nguyễn Văn Nam", "Đăng xuất"),
    @("20-10-2025 07:39:34", "Nguyễn Văn Nam", "Đăng nhập."),
    @("20-10-2025 07:39:49", "Nguyễn Văn Nam", "Đăng xuất"),
    @("20-10-2025 09:17:49", "Nguyễn Văn Nam", "Đăng nhập."),
    @("20-10-2025 09:18:13", "Nguyễn Văn Nam", "Đăng xuất"),
    @("20-10-2025 09:38:14", "Nguyễn Văn Nam", "Đăng nhập."),
    @("20-10-2025 09:39:39", "Nguyễn Văn Nam", "Đăng nhập."),
    @("20-10-2025 09:53:32", "Nguyễn Văn Nam", "Đăng nhập."),
    @("20-10-2025 09:53:45", "Nguyễn Văn Nam", "Thêm nhân viên: Đặng Thị Thúy (Staff)"),
    @("20-10-2025 09:53:53", "Nguyễn Văn Nam", "Quản lý xóa nhân viên: Đặng Thị Thúy"),
    @("20-10-2025 09:54:36", "Nguyễn Văn Nam", "Đăng nhập."),
    @("20-10-2025 09:54:51", "Nguyễn Văn Nam", "Thêm nhân viên: Phạm Thị Thúy (Staff)"),
    @("20-10-2025 09:54:55", "Nguyễn Văn Nam", "Đăng xuất"),
    @("20-10-2025 09:55:03", "Phạm Thị Thúy", "Đăng nhập.")
)

$startRow = 11
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
}

